$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2650
$ws.Range("I28").Value = 2475
$ws.Range("K28").Value = 2475
$ws.Range("M28").Value = -1990
$ws.Range("H113").Value = 4376.5835
$ws.Range("I113").Value = 4376.5835
$ws.Range("K113").Value = 4376.5835
$ws.Range("M113").Value = -1122.5835
$ws.Range("H115").Value = 996.3333
$ws.Range("I115").Value = 996.3333
$ws.Range("K115").Value = 2988.9999
$ws.Range("M115").Value = -1421.9999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1350.1794
$ws.Range("I32").Value = 1350.1794
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1350.1794
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = -1063.1794
$ws.Range("M32").ClearContents()
$ws.Range("H45").Value = 1923.5454
$ws.Range("I45").Value = 1917
$ws.Range("K45").Value = 1917
$ws.Range("M45").Value = -1540
$ws.Range("H122").Value = 4503.5
$ws.Range("I122").Value = 3905.2
$ws.Range("K122").Value = 11715.6
$ws.Range("M122").Value = -9265.599999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1579.2941
$ws.Range("I20").Value = 1181.1666
$ws.Range("K20").Value = 1181.1666
$ws.Range("M20").Value = -934.1666
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H97").Value = 22500
$ws.Range("J97").Value = 22500
$ws.Range("L97").Value = 22500
$ws.Range("N97").Value = -24482
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("N98").Value = 0
$ws.Range("L98").ClearContents()
$ws.Range("H99").Value = 2172.2307
$ws.Range("I99").Value = 2224.0833
$ws.Range("K99").Value = 2224.0833
$ws.Range("M99").Value = -726.0832999999998
$ws.Range("H100").Value = 15000
$ws.Range("J100").Value = 15000
$ws.Range("L100").Value = 15000
$ws.Range("N100").Value = -17164
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H103").Value = 87654
$ws.Range("J103").Value = 87654
$ws.Range("L103").Value = 87654
$ws.Range("N103").Value = -89998
$ws.Range("H105").Value = 3267.4614
$ws.Range("I105").Value = 3716
$ws.Range("J105").Value = 2883
$ws.Range("K105").Value = 3716
$ws.Range("L105").Value = 2883
$ws.Range("M105").Value = -1969
$ws.Range("N105").Value = -6377
$ws.Range("H106").Value = 50000
$ws.Range("J106").Value = 50000
$ws.Range("L106").Value = 50000
$ws.Range("N106").Value = -52524
$ws.Range("H107").Value = 37272.47
$ws.Range("I107").Value = 10124.615
$ws.Range("J107").Value = 125503
$ws.Range("K107").Value = 10124.615
$ws.Range("L107").Value = 125503
$ws.Range("M107").Value = -8204.615
$ws.Range("N107").Value = -129343

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 174.94737
$ws.Range("I7").Value = 118.9375
$ws.Range("K7").Value = 118.9375
$ws.Range("M7").Value = -5.9375
$ws.Range("H16").Value = 1377.7
$ws.Range("I16").Value = 1391.8889
$ws.Range("J16").Value = 1250
$ws.Range("K16").Value = 1391.8889
$ws.Range("L16").Value = 1250
$ws.Range("M16").Value = -1104.8889
$ws.Range("N16").Value = -1824
$ws.Range("H107").Value = 1070
$ws.Range("I107").Value = 875.125
$ws.Range("K107").Value = 875.125
$ws.Range("M107").Value = 1044.875
$ws.Range("H113").Value = 1377.7
$ws.Range("I113").Value = 1391.8889
$ws.Range("J113").Value = 1250
$ws.Range("K113").Value = 1391.8889
$ws.Range("L113").Value = 1250
$ws.Range("M113").Value = 778.1111000000001
$ws.Range("N113").Value = -5590
$ws.Range("H122").Value = 17734.467
$ws.Range("I122").Value = 1044.3334
$ws.Range("J122").Value = 84495
$ws.Range("K122").Value = 3133.0002
$ws.Range("L122").Value = 253485
$ws.Range("M122").Value = -683.0001999999999
$ws.Range("N122").Value = -258385

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 10700
$ws.Range("I76").Value = 4250
$ws.Range("J76").Value = 15000
$ws.Range("K76").Value = 12750
$ws.Range("L76").Value = 45000
$ws.Range("M76").Value = -12367
$ws.Range("N76").Value = -45766
$ws.Range("H79").Value = 10700
$ws.Range("I79").Value = 4250
$ws.Range("J79").Value = 15000
$ws.Range("K79").Value = 12750
$ws.Range("L79").Value = 45000
$ws.Range("M79").Value = -11424
$ws.Range("N79").Value = -47652
$ws.Range("H80").Value = 4750
$ws.Range("I80").Value = 4500
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 13500
$ws.Range("L80").Value = 15000
$ws.Range("M80").Value = -12564
$ws.Range("N80").Value = -16872
$ws.Range("H81").Value = 3563.1428
$ws.Range("I81").Value = 2657.1667
$ws.Range("K81").Value = 7971.500100000001
$ws.Range("M81").Value = -6848.500100000001
$ws.Range("H83").Value = 4750
$ws.Range("I83").Value = 4500
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 40500
$ws.Range("L83").Value = 45000
$ws.Range("M83").Value = -35820
$ws.Range("N83").Value = -54360
$ws.Range("H84").Value = 3563.1428
$ws.Range("I84").Value = 2657.1667
$ws.Range("K84").Value = 23914.5003
$ws.Range("M84").Value = -18298.5003
$ws.Range("H94").Value = 8536.375
$ws.Range("I94").Value = 7658.2
$ws.Range("K94").Value = 22974.6
$ws.Range("M94").Value = -22298.6
$ws.Range("H105").Value = 31472
$ws.Range("J105").Value = 31472
$ws.Range("L105").Value = 94416
$ws.Range("N105").Value = -99658
$ws.Range("H112").Value = 4855.5713
$ws.Range("H115").Value = 8875
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3338.0908
$ws.Range("I102").Value = 2790.087
$ws.Range("K102").Value = 2790.087
$ws.Range("M102").Value = -1168.087
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2018.4166
$ws.Range("I61").Value = 1561.9412
$ws.Range("K61").Value = 1561.9412
$ws.Range("M61").Value = -1359.9412
$ws.Range("H113").Value = 2018.4166
$ws.Range("I113").Value = 1561.9412
$ws.Range("K113").Value = 1561.9412
$ws.Range("M113").Value = 608.0588
$ws.Range("H132").Value = 3707
$ws.Range("I132").Value = 3370.2307
$ws.Range("K132").Value = 10110.6921
$ws.Range("M132").Value = -7580.6921

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 567
$ws.Range("I113").Value = 682.44446
$ws.Range("K113").Value = 2047.33338
$ws.Range("M113").Value = 122.66662
$ws.Range("H132").Value = 1719.0526
$ws.Range("I132").Value = 1770.1111
$ws.Range("K132").Value = 5310.3333
$ws.Range("M132").Value = -2780.3333
